$wb = $excel.ActiveWorkbook

# --- "名称" (Names) sheet: rename NPC C from "npc模板2" to "鱼线龙" ---
$wsNames = $wb.Worksheets.Item("名称")
$wsNames.Range("C4").Value = "鱼线龙"

# --- "描述" (Description) sheet: rename NPC C description to "申必人" ---
$wsDesc = $wb.Worksheets.Item("描述")
$wsDesc.Range("C4").Value = "申必人"

# --- "属性" (Attributes) sheet: new battle-system stat blocks for all three NPCs ---
$wsAttr = $wb.Worksheets.Item("属性")
$wsAttr.Range("B4").Value = "{""health"":{""躯干"":100,""左肢"":100,""右肢"":100,""左腿"":100,""右腿"":100,""头"":100},""法力"":1,'力量':22,'敏捷':0,'智力':0,'气运':100,'防御':0,'可以战斗':0}"
$wsAttr.Range("C4").Value = "{""health"":{""躯干"":1000,""左肢"":100,""右肢"":100,""左腿"":100,""右腿"":100,""头"":1000},""法力"":1,'力量':1,'敏捷':0,'智力':30,'气运':100,'防御':0,'可以战斗':0}"
$wsAttr.Range("D4").Value = "{""health"":{""躯干"":100,""左肢"":100,""右肢"":100,""左腿"":100,""右腿"":100,""头"":100},""法力"":1,'力量':19,'敏捷':0,'智力':0,'气运':100,'防御':20,'可以战斗':0}"
$wsAttr.Range("C4").Select()

# --- Add new "技能" (Skills) sheet at the end of the workbook ---
$wsFuncId = $wb.Worksheets.Item("功能ID")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsFuncId.Copy($null, $lastSheet)
$wsSkills = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsSkills.Name = "技能"
$wsSkills.Rows.Item(39).Delete()

$wsSkills.Range("B4").Value = "[""挥拳""]"
$wsSkills.Range("C4").Value = "[""挥拳"",""吟唱""]"
$wsSkills.Range("D4").Value = "[""挥拳""]"
$wsSkills.Range("B4").HorizontalAlignment = -4108
$wsSkills.Range("C4").Select()

# --- restore the originally active sheet / selection ---
$wsNames.Activate()
$wsNames.Range("D4").Select()
